# [Add]_02/SkillPopup 초상화, 아이콘 적용
# 1. SkillPopup 초상화, 아이콘 적용
#
# Adds an "IconPath" column (N) to the SkillData sheet with icon paths per
# skill row, and adjusts the EnemyData sheet's D2 value + both sheets'
# active-selection state (EnemyData becomes the active tab/sheet).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("SkillData")
$ws2 = $wb.Worksheets.Item("EnemyData")

# --- SkillData (sheet1): add IconPath column (N) -----------------------

# Carry over the same cell formatting/style that column M (EffectPath) uses
# for each row, so the new column N matches the look of its neighbour.
$ws1.Range("M1:M43").Copy()
$ws1.Range("N1:N43").PasteSpecial(-4122)  # xlPasteFormats

$ws1.Range("N1").Value = "IconPath"
$ws1.Range("N2").Value = "UI/Icon/Skill1"
$ws1.Range("N3").Value = "UI/Icon/Skill2"
$ws1.Range("N4").Value = "UI/Icon/Skill_Breesha3"
$ws1.Range("N5").Value = "UI/Icon/Skill_Breesha3"
$ws1.Range("N6").Value = "UI/Icon/Skill_Breesha3"
$ws1.Range("N7").Value = "UI/Icon/Skill_Breesha4"
$ws1.Range("N8").Value = "UI/Icon/Skill_Breesha4"
$ws1.Range("N9").Value = "UI/Icon/Skill1"
$ws1.Range("N10").Value = "UI/Icon/Skill2"
$ws1.Range("N11").Value = "UI/Icon/Skill_Eve3"
$ws1.Range("N12").Value = "UI/Icon/Skill_Eve3"
$ws1.Range("N13").Value = "UI/Icon/Skill_Eve3"
$ws1.Range("N14").Value = "UI/Icon/Skill_Eve4"
$ws1.Range("N15").Value = "UI/Icon/Skill_Eve4"
$ws1.Range("N16").Value = "UI/Icon/Skill1"
$ws1.Range("N17").Value = "UI/Icon/Skill2"
$ws1.Range("N18").Value = "UI/Icon/Skill_Adam3"
$ws1.Range("N19").Value = "UI/Icon/Skill_Adam3"
$ws1.Range("N20").Value = "UI/Icon/Skill_Adam3"
$ws1.Range("N21").Value = "UI/Icon/Skill_Adam4"
$ws1.Range("N22").Value = "UI/Icon/Skill_Adam4"
$ws1.Range("N23").Value = "UI/Icon/Skill1"
$ws1.Range("N24").Value = "UI/Icon/Skill2"
$ws1.Range("N25").Value = "UI/Icon/Skill_Abel3"
$ws1.Range("N26").Value = "UI/Icon/Skill_Abel3"
$ws1.Range("N27").Value = "UI/Icon/Skill_Abel3"
$ws1.Range("N28").Value = "UI/Icon/Skill_Abel4"
$ws1.Range("N29").Value = "UI/Icon/Skill_Abel4"
$ws1.Range("N30").Value = "UI/Icon/Skill1"
$ws1.Range("N31").Value = "UI/Icon/Skill2"
$ws1.Range("N32").Value = "UI/Icon/Skill_Kain3"
$ws1.Range("N33").Value = "UI/Icon/Skill_Kain3"
$ws1.Range("N34").Value = "UI/Icon/Skill_Kain3"
$ws1.Range("N35").Value = "UI/Icon/Skill_Kain4"
$ws1.Range("N36").Value = "UI/Icon/Skill_Kain4"
$ws1.Range("N37").Value = "UI/Icon/Skill1"
$ws1.Range("N38").Value = "UI/Icon/Skill2"
# Rows 39-43 intentionally left blank (no IconPath for those skill rows),
# matching the formatting-only cells carried over above.

# Widen the new column to fit its longest value.
$ws1.Columns.Item(14).ColumnWidth = 25.17

# --- EnemyData (sheet2): tweak a data value -----------------------------

$ws2.Range("D2").Value = 200

# --- Selection / active-tab state --------------------------------------
# SkillData keeps a specific selection, EnemyData becomes the active sheet.

$ws1.Range("J14").Select()
$ws2.Activate()
$ws2.Range("H16").Select()
